# UC009 - Prestar Contas (GT) — v1.2.1 -> v1.2.3
#
# TC3/TC4: swap the 3rd step's Action/Expected Result pair so that
# TC3 ends with "detalhar a solicitação de diária" and TC4 ends with
# "excluir comprovante" (they were previously reversed).
#
# TC7/TC8/TC9: rotate the "Expected Results" text so that TC7 gets the
# message that used to belong to TC9, and TC9 gets the message that
# used to belong to TC7 (TC8's expected result is unaffected).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# TC3 - step 3 (row 32): was "excluir comprovante" -> now "detalhar a solicitação de diária"
$ws.Range("B32").Value = "Chefe Clica para detalhar a solicitação de diária."
$ws.Range("D32").Value = "SYSTEM Apresenta a tela de Detalhar Diárias"

# TC4 - step 3 (row 41): was "detalhar a solicitação de diária" -> now "excluir comprovante"
$ws.Range("B41").Value = "Chefe Clica em excluir comprovante."
$ws.Range("D41").Value = "SYSTEM Exclui o comprovante."

# TC7 - step 2 expected result (row 67): now gets TC9's old message
$ws.Range("D67").Value = "SYSTEM Identifica que a prestação de contas indicada pelo usuário não está em nenhum desses dois estados: a) NÃO REALIZADA e b) DEVOLVIDA; Permite não permite um novo envio ou alterações na prestação (exclusão de documentos)."

# TC9 - step 2 expected result (row 82): now gets TC7's old message
$ws.Range("D82").Value = "SYSTEM Identifica que a solicitação indicada pelo usuário ainda não pode ter sua prestação de contas realizada; Exibe mensagem de erro (MSG212 - Prestação de contas ainda não pode ser realizada) para o usuário, impedindo que ele preste contas (anexa arquivos e etc)."
